# Automatic update of files.
# Bump the "Förändrad" (Changed) date in column C for every data row
# from 2023-11-03 (serial 45233) to 2023-11-13 (serial 45243).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value = 45243
    }
}
